$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename to snake_case machine-friendly column names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Normalize capitalization of connector words (de/el/la/los -> De/El/La/Los) ---
$ws.Range("A17").Value = "Ciudad De México"
$ws.Range("A28").Value = "Estado De México"
$ws.Range("B29").Value = "Atizapán De Zaragoza"
$ws.Range("B30").Value = "Ixtapan De La Sal"
$ws.Range("B38").Value = "Apaseo El Alto"
$ws.Range("B42").Value = "Acapulco De Juárez"
$ws.Range("B46").Value = "Ayutla De Los Libres"
$ws.Range("B49").Value = "Huitzuco De Los Figueroa"
$ws.Range("B50").Value = "Iguala De La Independencia"
$ws.Range("B51").Value = "Zihuatanejo De Azueta"
$ws.Range("B56").Value = "Tepecoacuilco De Trujano"
$ws.Range("B57").Value = "Tlapa De Comonfort"
$ws.Range("B61").Value = "Atotonilco De Tula"
$ws.Range("B67").Value = "Nopala De Villagrán"
$ws.Range("B69").Value = "Tepehuacán De Guerrero"
$ws.Range("B77").Value = "Tizapán El Alto"
$ws.Range("B78").Value = "Unión De Tula"
$ws.Range("B97").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B104").Value = "Cuayuca De Andrade"
$ws.Range("B107").Value = "Izúcar De Matamoros"
$ws.Range("B110").Value = "San Nicolás De Los Ranchos"
$ws.Range("B116").Value = "Landa De Matamoros"
$ws.Range("B138").Value = "Ignacio De La Llave"

# --- Remove trailing footer/metadata rows (153-157) ---
$ws.Range("A153:A157").EntireRow.Delete()
